$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "25.751.68"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.628.63"
$ws.Range("E3").Value = "  +0.15%  "
Set-TextValue "D4" "0.997"
Set-TextValue "D5" "213.83"
$ws.Range("E5").Value = "  -0.18%  "
Set-TextValue "D6" "0.501"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("E8").Value = "  -0.87%  "
Set-TextValue "D9" "0.0631"
$ws.Range("E9").Value = "  -0.54%  "
Set-TextValue "D10" "19.63"
$ws.Range("E10").Value = "  +0.47%  "
Set-TextValue "D11" "0.0789"
$ws.Range("E11").Value = "  +0.74%  "
Set-TextValue "D12" "4.24"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "1.853.57"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "1.629.47"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "0.0₃0759"
$ws.Range("E16").Value = "  -0.23%  "
Set-TextValue "D17" "62.68"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "25.748.91"
$ws.Range("E18").Value = "  +0.26%  "
Set-TextValue "D19" "0.997"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("E20").Value = "  +0.33%  "
Set-TextValue "D21" "191.05"
$ws.Range("E21").Value = "  -1.05%  "
Set-TextValue "D22" "9.91"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("E25").Value = "  +1.34%  "
Set-TextValue "D26" "142.07"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("E28").Value = "  +0.07%  "
Set-TextValue "D29" "15.50"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("E30").Value = "  +0.00%  "
Set-TextValue "D31" "0.0494"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("E32").Value = "  +0.20%  "
Set-TextValue "D33" "3.21"
$ws.Range("E33").Value = "  -0.50%  "
Set-TextValue "D34" "1.59"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  -0.46%  "
Set-TextValue "D36" "0.902"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("D37").Value = "1.142.25"
$ws.Range("E37").Value = "  +3.34%  "
Set-TextValue "D38" "0.544"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  -2.10%  "
Set-TextValue "D40" "0.0155"
$ws.Range("E40").Value = "  +0.15%  "
Set-TextValue "D41" "0.996"
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("E43").Value = "  +0.88%  "
Set-TextValue "D44" "101.01"
$ws.Range("E44").Value = "  +1.10%  "
Set-TextValue "D45" "0.800"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "1.763.84"
$ws.Range("E46").Value = "  +0.35%  "
Set-TextValue "D47" "55.17"
$ws.Range("E47").Value = "  +0.58%  "
Set-TextValue "D48" "1.47"
$ws.Range("E48").Value = "  +7.68%  "
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "7.51"
$ws.Range("E51").Value = "  -1.86%  "
